$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# CONFIG sheet: "defaults.state_fallback" (row 11) never had a value in the
# value column (B). Touch it with a formatting no-op so the (empty) cell is
# materialized without altering its style.
# ---------------------------------------------------------------------------
$wsConfig = $wb.Worksheets.Item("CONFIG")
$wsConfig.Cells.Item(11, 2).Font.Bold = $false

# ---------------------------------------------------------------------------
# INPUT_MASTER sheet: make room for 17 new lead/ECorp columns in front of the
# existing BatchData (BD_*) columns, then populate the new headers and carry
# over the bold/bordered/centered header style used by the existing header
# row.
# ---------------------------------------------------------------------------
$wsInput = $wb.Worksheets.Item("INPUT_MASTER")

# Shift the current A1:T1 header block 17 columns to the right (-> R1:AK1),
# opening up A1:Q1 for the new headers.
$wsInput.Range("A1:Q1").Insert(-4161)

$newHeaders = @(
    "FULL_ADDRESS",
    "COUNTY",
    "Owner_Ownership",
    "ECORP_INDEX_#",
    "OWNER_TYPE",
    "ECORP_SEARCH_NAME",
    "ECORP_TYPE",
    "ECORP_NAME_S",
    "ECORP_ENTITY_ID_S",
    "ECORP_ENTITY_TYPE",
    "ECORP_STATUS",
    "ECORP_FORMATION_DATE",
    "ECORP_BUSINESS_TYPE",
    "ECORP_STATE",
    "ECORP_COUNTY",
    "ECORP_COMMENTS",
    "ECORP_URL"
)

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $wsInput.Cells.Item(1, $i + 1).Value = $newHeaders[$i]
}

# Clone the header formatting (bold font, thin border, centered) from the
# existing header cells onto the newly-inserted header cells.
$wsInput.Range("R1").Copy()
$wsInput.Range("A1:Q1").PasteSpecial(-4122)
